$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: "15455" -> "15092" (keep stored as text, matching original inlineStr type)
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "15092"
$ws.Range("B2").Style = "Normal"

# C2: ATA code
$ws.Range("C2").Value = "21-61"

# D2: B1-code
$ws.Range("D2").Value = "B1-005812"

# E2: LRU
$ws.Range("E2").Value = "ACSC 1"

# F2: Message
$ws.Range("F2").Value = "INTERNAL FAULT"

# H2: Potential FDE -> single space
$ws.Range("H2").Value = " "

# I2: Date From changes, Date To (J2) unchanged
$ws.Range("I2").Value = 44161

# L2: ISE Input
$l2 = @"
Per SL-21-018 (to be revised as of 9/24/20)
This is related to pack cycling. Will have to revise SL procedure (by 30 Sept 2020) and the FIM (31 Oct 2020). (Input from specialist Sep/2020)
Fleet wide msg in top #10 position
Reset SL procedure is not in FIM
"@
$ws.Range("L2").Value = $l2

# M2: ISE Rec Act
$m2 = @"
1. Reset per SL procedure: 
a) Reset procedure
b) Wait for 30 sec., then select L Pack Manual Mode 
c) Wait for 30 sec., and then select Pack Auto Mode. 
d) If still cycling Swap ACSC. If not replace ACSC. Do not remove sensors
2. Follow FIM 21−61−04−810−81:
NOTE:Bombardier strongly recommends to swap ASCSs before replacing the LRUs. Byswapping LRUs, you will verify the integrity of the ASCSs. Go to Fault Confirmation.
1. Pack Discharge Pressure Sensor (PDPS) MT13 unserviceable.
2. Pack Inlet Flow Sensor (PIFS) MT11 unserviceable.
3. Pack Inlet Pressure Sensor (PIPS) MT9 unserviceable.
4. Defective wiring interface between Air Conditioning System Controller #1 (ACSC 1)and oneof the pressure sensors.
"@
$ws.Range("M2").Value = $m2
